$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
Write-Host $ws.Name
